# daily auto push: 2026-02-21 04:12 UTC
# A new measurement row for 2026/02/21 (土) is inserted right after the
# existing 2026/02/21 rows (row 825), pushing every subsequent row down
# by one (826 -> 827, ..., 867 -> 868) and extending the used range from
# D867 to D868.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the last existing "2026/02/21" row (825) and use "insert copied
# cells" semantics so the new row 826 inherits the same date/weekday text
# formatting (kept as literal text, not auto-converted to a date serial)
# and no extra styles are introduced.
$ws.Rows("825:825").Copy()
$ws.Rows("826:826").Insert()

# Overwrite the time/rank columns for the newly inserted row; column A
# ("2026/02/21") and column B ("土") already carry the correct values from
# the copy, so only C and D need to change.
$ws.Range("C826").Value = 12
$ws.Range("D826").Value = 88

$excel.CutCopyMode = $false
